{"js": "// Apply the PMO Milestone Schedule edits described by the commit/diff:\n//  1. Milestone 1: collapse the split \"Solution set up with [Xen] Engine & HGE\n//     Engine.\" runs (which wrap \"Xen\" in proofErr spell-check markup) into a\n//     single run, keeping the existing strikethrough formatting.\n//  2. Milestone 6: insert two new bullets (\"Polish Ghost and PacMan AI.\" then\n//     \"Add metadata encryption.\") right after \"Finalise all sound effects and\n//     music.\" and before \"Release the first revision with STANDARD and PRO\n//     versions.\", and remove the old, now-duplicate \"Polish Ghost and PacMan\n//     AI.\" bullet (plus its trailing blank paragraph) that used to sit just\n//     before the \"Milestone 7\" heading.\n//  3. Milestone 7: retitle the last bullet from \"Full schedule from this\n//     point TBD.\" to \"Add in the tutorial screens to the menu.\"\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst texts = paragraphs.items.map((p) => p.text);\n\nfunction findIndex(target, fromIndex) {\n  for (let i = fromIndex || 0; i < texts.length; i++) {\n    if (texts[i].trim() === target) {\n      return i;\n    }\n  }\n  throw new Error(\"Paragraph not found: \" + target);\n}\n\n// --- Change 1: merge the \"Solution set up with Xen Engine & HGE Engine.\"\n// runs (and drop the spell-check proofErr markup around \"Xen\") into one run.\nconst solutionIdx = findIndex(\"Solution set up with Xen Engine & HGE Engine.\");\nparagraphs.items[solutionIdx]\n  .getRange()\n  .insertText(\"Solution set up with Xen Engine & HGE Engine.\", \"Replace\");\nawait context.sync();\n\n// --- Change 2a: add the two new bullets after \"Finalise all sound effects\n// and music.\" (still inside Milestone 6's numId=5 list).\nconst finaliseIdx = findIndex(\"Finalise all sound effects and music.\");\nconst ghostAiPara = paragraphs.items[finaliseIdx].insertParagraph(\n  \"Polish Ghost and PacMan AI.\",\n  \"After\"\n);\nawait context.sync();\nghostAiPara.insertParagraph(\"Add metadata encryption.\", \"After\");\nawait context.sync();\n\n// --- Change 2b: remove the old \"Polish Ghost and PacMan AI.\" bullet (and the\n// blank paragraph right after it) that used to precede the \"Milestone 7\"\n// heading.\nparagraphs.load(\"text\");\nawait context.sync();\nconst refreshedTexts = paragraphs.items.map((p) => p.text);\nlet oldGhostAiIdx = -1;\nfor (let i = 0; i < refreshedTexts.length; i++) {\n  if (refreshedTexts[i].trim() === \"Polish Ghost and PacMan AI.\") {\n    oldGhostAiIdx = i;\n  }\n}\nif (oldGhostAiIdx === -1) {\n  throw new Error(\"Could not find the old 'Polish Ghost and PacMan AI.' bullet\");\n}\nconst blankAfter = paragraphs.items[oldGhostAiIdx + 1];\nblankAfter.load(\"text\");\nawait context.sync();\nif (blankAfter.text.trim() !== \"\") {\n  throw new Error(\"Expected a blank paragraph after the old bullet to remove\");\n}\nblankAfter.delete();\nparagraphs.items[oldGhostAiIdx].delete();\nawait context.sync();\n\n// --- Change 3: rename the final Milestone 7 bullet.\nparagraphs.load(\"text\");\nawait context.sync();\nconst tbdIdx = findIndex(\"Full schedule from this point TBD.\");\nparagraphs.items[tbdIdx]\n  .getRange()\n  .insertText(\"Add in the tutorial screens to the menu.\", \"Replace\");\nawait context.sync();\n", "ps1": "# Apply the PMO Milestone Schedule edits described by the commit/diff:\n#  1. Milestone 1: collapse the split \"Solution set up with [Xen] Engine & HGE\n#     Engine.\" runs (which wrap \"Xen\" in proofErr spell-check markup) into a\n#     single run, keeping the existing strikethrough formatting.\n#  2. Milestone 6: insert two new bullets (\"Polish Ghost and PacMan AI.\" then\n#     \"Add metadata encryption.\") right after \"Finalise all sound effects and\n#     music.\" and before \"Release the first revision with STANDARD and PRO\n#     versions.\", and remove the old, now-duplicate \"Polish Ghost and PacMan\n#     AI.\" bullet (plus its trailing blank paragraph) that used to sit just\n#     before the \"Milestone 7\" heading.\n#  3. Milestone 7: retitle the last bullet from \"Full schedule from this\n#     point TBD.\" to \"Add in the tutorial screens to the menu.\"\n\n$d = $word.ActiveDocument\n\nfunction Find-ParagraphIndex($doc, $targetText) {\n    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {\n        if ($doc.Paragraphs.Item($i).Range.Text.Trim() -eq $targetText) {\n            return $i\n        }\n    }\n    throw \"Paragraph not found: $targetText\"\n}\n\n# --- Change 1: merge the \"Solution set up with Xen Engine & HGE Engine.\"\n# runs (and drop the spell-check proofErr markup around \"Xen\") into one run.\n$solutionIdx = Find-ParagraphIndex $d \"Solution set up with Xen Engine & HGE Engine.\"\n$solutionRange = $d.Paragraphs.Item($solutionIdx).Range\n$find = $solutionRange.Find\n$find.Execute(\"Solution set up with Xen Engine & HGE Engine.\", $false, $false, $false, $false, $false, $true, 1, $false, \"Solution set up with Xen Engine & HGE Engine.\", 2) | Out-Null\n\n# --- Change 2a: add the two new bullets after \"Finalise all sound effects\n# and music.\" (still inside Milestone 6's numId=5 list).\n$finaliseIdx = Find-ParagraphIndex $d \"Finalise all sound effects and music.\"\n$d.Paragraphs.Item($finaliseIdx).Range.InsertParagraphAfter()\n$ghostAiPara = $d.Paragraphs.Item($finaliseIdx + 1)\n$ghostAiPara.Range.Text = \"Polish Ghost and PacMan AI.\"\n$ghostAiPara.Range.InsertParagraphAfter()\n$metadataPara = $d.Paragraphs.Item($finaliseIdx + 2)\n$metadataPara.Range.Text = \"Add metadata encryption.\"\n\n# --- Change 2b: remove the old \"Polish Ghost and PacMan AI.\" bullet (and the\n# blank paragraph right after it) that used to precede the \"Milestone 7\"\n# heading. (Search from after the newly-inserted bullet so we find the old\n# occurrence, not the one we just added.)\n$oldGhostAiIdx = -1\nfor ($i = $finaliseIdx + 3; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text.Trim() -eq \"Polish Ghost and PacMan AI.\") {\n        $oldGhostAiIdx = $i\n        break\n    }\n}\nif ($oldGhostAiIdx -eq -1) {\n    throw \"Could not find the old 'Polish Ghost and PacMan AI.' bullet\"\n}\n$blankText = $d.Paragraphs.Item($oldGhostAiIdx + 1).Range.Text.Trim()\nif ($blankText -ne \"\") {\n    throw \"Expected a blank paragraph after the old bullet to remove\"\n}\n$d.Paragraphs.Item($oldGhostAiIdx + 1).Range.Delete() | Out-Null\n$d.Paragraphs.Item($oldGhostAiIdx).Range.Delete() | Out-Null\n\n# --- Change 3: rename the final Milestone 7 bullet.\n$tbdIdx = Find-ParagraphIndex $d \"Full schedule from this point TBD.\"\n$tbdRange = $d.Paragraphs.Item($tbdIdx).Range\n$find2 = $tbdRange.Find\n$find2.Execute(\"Full schedule from this point TBD.\", $false, $false, $false, $false, $false, $true, 1, $false, \"Add in the tutorial screens to the menu.\", 2) | Out-Null\n"}
